$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.613414406776428
$ws.Range("B1").Value = 3.298765897750854
$ws.Range("C1").Value = 4.273428916931152
$ws.Range("D1").Value = 1.328328132629395
$ws.Range("E1").Value = 0.7779595255851746
